# Fills in the "English" and "Math" subject tables' comment and grade
# cells in the student's report card. The document has 4 subject tables
# (Torah, English, Math, Geography), each with the same 2x2 layout:
#   Cell(1,1)=subject name   Cell(1,2)=teacher comment (merged, 2 rows)
#   Cell(2,1)=["ציון:" paragraph, grade paragraph]   Cell(2,2)=merged w/ (1,2)
#
# NOTE ON RUNTIME QUIRK: in this COM-interop runtime, once $d.Tables (or any
# Table/Cell) has been accessed, later calls to $d.Paragraphs.Item(n) return
# stale/incorrect ranges. So step 1 below resolves the "grade" paragraph
# indexes purely via the document-level Paragraphs collection, and the
# grades are written *before* Tables is touched at all. Table objects are
# also re-fetched fresh via $d.Tables.Item(i) rather than cached in
# variables that stay alive across statements, to avoid the same issue.

$d = $word.ActiveDocument

# --- STEP 1: locate the "grade" paragraph for every subject table. It is
# the (empty) paragraph that immediately follows the paragraph beginning
# with "ציון:" inside Cell(2,1) of each table.
$tzionChar = [string][char]1510  # Hebrew Tsadi, first letter of "ציון"
$gradeParaIndex = New-Object System.Collections.ArrayList
$prevWasTzion = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($prevWasTzion -and $t.Length -le 2) {
        [void]$gradeParaIndex.Add($i)
    }
    $prevWasTzion = ($t.Length -ge 1) -and ($t.Substring(0, 1) -eq $tzionChar)
}

# Tables/grades appear in document order: Torah, English, Math, Geography
$gradeEnglishIdx = $gradeParaIndex[1]
$gradeMathIdx = $gradeParaIndex[2]

# --- STEP 2: write the two grades (must happen before any Tables access).
$d.Paragraphs.Item($gradeEnglishIdx).Range.Text = "85"
$d.Paragraphs.Item($gradeMathIdx).Range.Text = "92"

# --- STEP 3: write the two teacher comments, matching tables by subject name.
for ($i = 1; $i -le $d.Tables.Count; $i++) {
    $subject = $d.Tables.Item($i).Cell(1,1).Range.Text
    if ($subject.StartsWith("אנגלית")) {
        $d.Tables.Item($i).Cell(1,2).Range.Text = "במחצית זאת חזרנו על הדקדוק, עבר הווה ועתיד, חזרנו על השיטות שפיתחנו להבנת הניקרא, עברנו על מאמרים קשים, וניסנו לפתור אותם בשיטות לימוד, תרגלנו הרבה לקראת מבחני גמר!`nריקי את מעולה, הרבה הצלחה "
    } elseif ($subject.StartsWith("חשבון")) {
        $d.Tables.Item($i).Cell(1,2).Range.Text = "במחצית זאת למדנו משפט פיתגורס,והרחבנו בנושא חפיפת משולשים, חזרנו על משוואות ב2 נעלמים, ניתוח גרפים והתכוננו לקראת מבחני גמר במתמטיקה.`nריקי את מעולה, בהצלחה רבה בהמשך!"
    }
}
